$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 41
$ws.Range("C2").Value = 2
$ws.Range("E2").Formula = "=36/60"
$ws.Range("E2").HorizontalAlignment = -4152
$ws.Rows.Item(2).RowHeight = 13.75
